$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-43) holds the "Förändrad" date, stored as date serial 45767
# (2025-04-20). Bump it by one day to 45768 (2025-04-21) for every data row.
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45767) {
        $cell.Value2 = 45768
    }
}
